$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly refresh of the cryptos table (GitHub Actions cron).
# Column D ("Price") holds numeric-looking text (e.g. "1.007", "28.259.59")
# that must stay text, not be coerced to a real number -- so those writes are
# apostrophe-prefixed exactly like typing text into Excel by hand.
# Column E ("Volume(1h)") is already unambiguous text (%, spaces) so no prefix
# is needed there.

$ws.Range("D2").Value = "'28.259.59"
$ws.Range("E2").Value = '  +0.88%  '

$ws.Range("D3").Value = "'1.884.10"

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = "'314.06"
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").Value = "'0.5142"
$ws.Range("E7").Value = '  +1.00%  '

$ws.Range("D8").Value = "'0.3907"
$ws.Range("E8").Value = '  +2.82%  '

$ws.Range("D9").Value = "'0.08377"
$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("E10").Value = '  +0.98%  '

$ws.Range("D11").Value = "'41.67"
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").Value = "'6.241"
$ws.Range("E12").Value = '  +0.41%  '

$ws.Range("D13").Value = "'20.76"
$ws.Range("E13").Value = '  +1.50%  '

$ws.Range("D14").Value = "'1.882.58"
$ws.Range("E14").Value = '  +0.74%  '

$ws.Range("D15").Value = "'7.294"
$ws.Range("E15").Value = '  +1.62%  '

$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = '  +0.32%  '

$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("E18").Value = '  +0.93%  '

$ws.Range("D19").Value = "'0.06667"
$ws.Range("E19").Value = '  +0.69%  '

$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").Value = "'1.006"
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").Value = "'6.078"
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").Value = "'28.293.46"
$ws.Range("E23").Value = '  +0.86%  '

$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = '  +0.78%  '

$ws.Range("D26").Value = "'2.094.10"
$ws.Range("E26").Value = '  +0.46%  '

$ws.Range("E27").Value = '  -1.80%  '

$ws.Range("D28").Value = "'158.89"
$ws.Range("E28").Value = '  +1.03%  '

$ws.Range("E29").Value = '  +1.07%  '

$ws.Range("D30").Value = "'125.61"
$ws.Range("E30").Value = '  +0.09%  '

$ws.Range("D31").Value = "'0.1067"
$ws.Range("E31").Value = '  +1.06%  '

$ws.Range("E32").Value = '  +0.71%  '

$ws.Range("E33").Value = '  +5.49%  '

$ws.Range("D34").Value = "'3.597"
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = "'9.782"
$ws.Range("E35").Value = '  +1.21%  '

$ws.Range("D36").Value = "'0.02454"
$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("D37").Value = "'0.06583"
$ws.Range("E37").Value = '  +0.83%  '

$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("D39").Value = "'1.214"
$ws.Range("E39").Value = '  +0.56%  '

$ws.Range("D40").Value = "'0.6549"
$ws.Range("E40").Value = '  +2.12%  '

$ws.Range("D41").Value = "'5.037"
$ws.Range("E41").Value = '  +3.51%  '

$ws.Range("D42").Value = "'1.234"
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("D43").Value = "'11.33"
$ws.Range("E43").Value = '  +0.71%  '

$ws.Range("D44").Value = "'0.6144"
$ws.Range("E44").Value = '  +0.74%  '

$ws.Range("D45").Value = "'13.19"
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("D46").Value = "'1.289"
$ws.Range("E46").Value = '  +0.39%  '

$ws.Range("D47").Value = "'3.679"
$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("D48").Value = "'2.018"
$ws.Range("E48").Value = '  +1.09%  '

$ws.Range("E49").Value = '  +2.21%  '

$ws.Range("D50").Value = "'121.76"
$ws.Range("E50").Value = '  +0.35%  '

$ws.Range("D51").Value = "'79.02"
$ws.Range("E51").Value = '  -0.85%  '
